$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing "Execution Time" timestamps (re-run finished later) ---
$ws.Cells.Item(2,5).Value = "28/03/2025 03:46:14 PM"
$ws.Cells.Item(3,5).Value = "28/03/2025 03:46:17 PM"
$ws.Cells.Item(4,5).Value = "28/03/2025 03:46:18 PM"
$ws.Cells.Item(5,5).Value = "28/03/2025 03:46:21 PM"

# --- Row 6 used to be the "Click Go to Sign In" sub-step; it is now replaced
#     by a new sub-step describing the account lockout after repeated wrong
#     OTP entries. ---
$ws.Cells.Item(6,3).Value = "Account Blocked After Multiple Wrong OTPs"
$ws.Cells.Item(6,5).Value = "28/03/2025 03:46:25 PM"
$ws.Cells.Item(6,6).Value = "You have reached the maximum login attempts for the day. Please try again after 24 hours."

# --- The old "Click Go to Sign In" sub-step now shifts down into a brand new
#     row 7 (pushing the rest of the rows down by one). ---
$ws.Cells.Item(7,1).Value = 2
# Leading-apostrophe trick: forces an actual (non-null) empty text value,
# matching the blank "Test Case Name" sub-step cell used by the other rows.
$ws.Cells.Item(7,2).Formula = "'"
$ws.Cells.Item(7,3).Value = "Click Go to Sign In"
$ws.Cells.Item(7,4).Value = "PASSED"
$ws.Cells.Item(7,5).Value = "28/03/2025 03:46:25 PM"
$ws.Cells.Item(7,6).Value = "Navigated back to Get OTP page successfully"

# --- Test case 3 ("User Login with Invalid Credentials") now lives on row 8
#     and failed with a stale-element Selenium exception instead of passing. ---
$staleElementError = @'
stale element reference: stale element not found in the current frame
  (Session info: chrome=134.0.6998.166)
For documentation on this error, please visit: https://www.selenium.dev/documentation/webdriver/troubleshooting/errors#stale-element-reference-exception
Build info: version: '4.29.0', revision: '5fc1ec94cb'
System info: os.name: 'Windows 10', os.arch: 'amd64', os.version: '10.0', java.version: '23.0.1'
Driver info: org.openqa.selenium.chrome.ChromeDriver
Command: [f0be20835ff80e2a170a0b82c403494d, clickElement {id=f.00095D8C67BDCF3ABF05023232D37247.d.F5A6A02415C47404393A2FB587A4C4BE.e.31}]
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 134.0.6998.166, chrome: {chromedriverVersion: 134.0.6998.165 (fd886e2cb29..., userDataDir: C:\Users\Admin\AppData\Loca...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:58998}, networkConnectionEnabled: false, pageLoadStrategy: normal, platformName: windows, proxy: Proxy(), se:cdp: ws://localhost:58998/devtoo..., se:cdpVersion: 134.0.6998.166, setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Element: [[ChromeDriver: chrome on windows (f0be20835ff80e2a170a0b82c403494d)] -> xpath: //button[@id='loginButton']]
Session ID: f0be20835ff80e2a170a0b82c403494d
'@

$ws.Cells.Item(8,1).Value = 3
$ws.Cells.Item(8,2).Value = "User Login with Invalid Credentials"
$ws.Cells.Item(8,3).Value = "No Sub-Steps"
$ws.Cells.Item(8,4).Value = "FAILED"
$ws.Cells.Item(8,5).Value = "28/03/2025 03:46:26 PM"
$ws.Cells.Item(8,6).Value = $staleElementError

# --- Test case 4 ("Verify Home Page Loads Successfully") shifts down to the
#     new row 9. ---
$ws.Cells.Item(9,1).Value = 4
$ws.Cells.Item(9,2).Value = "Verify Home Page Loads Successfully"
$ws.Cells.Item(9,3).Value = "No Sub-Steps"
$ws.Cells.Item(9,4).Value = "PASSED"
$ws.Cells.Item(9,5).Value = "28/03/2025 03:46:26 PM"
$ws.Cells.Item(9,6).Value = "Test executed successfully."
